$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Seed rows 5-10 with the same base look as the existing imported rows
#    (row 4 uses style index 1 on every column) before touching any values,
#    so the new rows keep an explicit style like the rest of the import
#    instead of falling back to the untouched default style.
# ---------------------------------------------------------------------------
$ws.Range("A4:L4").Copy()
$ws.Range("A5:L10").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Write all new row data (rows 5-10), row-major / column A..L order, so
#    that shared-string indices are allocated in the same order as the
#    target workbook (new unique strings appear in exactly that sequence).
# ---------------------------------------------------------------------------

# Row 5
$ws.Range("A5").Value = "Calle 79B No. 7-60"
$ws.Range("B5").Value = "Calle de los anticuarios"
$ws.Range("C5").Value = 110221
$ws.Range("D5").Value = 524
$ws.Range("E5").Value = 11
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3213723648
$ws.Range("H5").Value = "Daniela Rivera"
$ws.Range("I5").Value = 4.661714
$ws.Range("J5").Value = -74.050852
$ws.Range("K5").Value = "shipping"
$ws.Range("L5").Value = 13

# Row 6
$ws.Range("A6").Value = "Calle 59 # 8 - 28"
$ws.Range("B6").Value = "Chapinero Central"
$ws.Range("C6").Value = 110221
$ws.Range("D6").Value = 524
$ws.Range("E6").Value = 11
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3022367474
$ws.Range("H6").Value = "John Angel"
$ws.Range("I6").Value = 4.644743
$ws.Range("J6").Value = -74.062218
$ws.Range("K6").Value = "shipping"
$ws.Range("L6").Value = 9

# Row 7
$ws.Range("A7").Value = "Cra 7 # 45 - 49"
$ws.Range("B7").Value = "Chapinero Central"
$ws.Range("C7").Value = 110221
$ws.Range("D7").Value = 524
$ws.Range("E7").Value = 11
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3126608380
$ws.Range("H7").Value = "John Angel"
$ws.Range("I7").Value = 4.632111
$ws.Range("J7").Value = -74.06397
$ws.Range("K7").Value = "shipping"
$ws.Range("L7").Value = 9

# Row 8
$ws.Range("A8").Value = "Carrera 1 #67-21"
$ws.Range("B8").Value = "edificio nueva granada etapa1 apto 102"
$ws.Range("C8").Value = 110221
$ws.Range("D8").Value = 524
$ws.Range("E8").Value = 11
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3108546364
$ws.Range("H8").Value = "Tatiana Luna"
$ws.Range("I8").Value = 4.645672
$ws.Range("J8").Value = -74.051827
$ws.Range("K8").Value = "shipping"
$ws.Range("L8").Value = 11

# Row 9
$ws.Range("A9").Value = "Calle 10 #11A- 25 sur"
$ws.Range("B9").Value = "Ciudad Berna"
$ws.Range("C9").Value = 110221
$ws.Range("D9").Value = 524
$ws.Range("E9").Value = 11
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3115243150
$ws.Range("H9").Value = "Jennifer Lasprilla"
$ws.Range("I9").Value = 4.584875
$ws.Range("J9").Value = -74.090762
$ws.Range("K9").Value = "shipping"
$ws.Range("L9").Value = 8

# Row 10
$ws.Range("A10").Value = "Calle 6D # 3-89"
$ws.Range("B10").Value = "Candelaria"
$ws.Range("C10").Value = 110221
$ws.Range("D10").Value = 524
$ws.Range("E10").Value = 11
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = "310 854 6364"
$ws.Range("H10").Value = "Mercedes Garcia"
$ws.Range("I10").Value = 4.592404
$ws.Range("J10").Value = -74.074813
$ws.Range("K10").Value = "shipping"
$ws.Range("L10").Value = 10

# ---------------------------------------------------------------------------
# 3. Column A width.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21.14

# ---------------------------------------------------------------------------
# 4. Formatting for the imported rows. Each new look is set once on a
#    "primary" cell (creating the new font/fill), then copied with
#    paste-special-formats onto the other cells sharing that look so the
#    style is reused instead of re-derived.
# ---------------------------------------------------------------------------

# Style A: G5 -> Oswald 11pt #777777 on white fill, left aligned.
$ws.Range("G5").Font.Name = "Oswald"
$ws.Range("G5").Font.Size = 11
$ws.Range("G5").Font.Color = 7829367
$ws.Range("G5").Interior.Color = 16777215
$ws.Range("G5").HorizontalAlignment = -4131

# Style B: A9 -> Roboto #000000 on white fill, general alignment.
# Reused (via copy/paste-format) on G9, A10, G10, H10.
$ws.Range("A9").Font.Name = "Roboto"
$ws.Range("A9").Font.Color = 0
$ws.Range("A9").Interior.Color = 16777215

$ws.Range("A9").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("G10").PasteSpecial(-4122)
$ws.Range("H10").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# paste-format only touches formatting, but re-assert values/text defensively
$ws.Range("G9").Value = 3115243150
$ws.Range("A10").Value = "Calle 6D # 3-89"
$ws.Range("G10").Value = "310 854 6364"
$ws.Range("H10").Value = "Mercedes Garcia"

# Style C: B9 -> Arial #000000 on white fill, left aligned.
$ws.Range("B9").Font.Name = "Arial"
$ws.Range("B9").Font.Color = 0
$ws.Range("B9").Interior.Color = 16777215
$ws.Range("B9").HorizontalAlignment = -4131

# Style D: B10 -> Arial 11pt #000000 on white fill, centered.
$ws.Range("B10").Font.Name = "Arial"
$ws.Range("B10").Font.Size = 11
$ws.Range("B10").Font.Color = 0
$ws.Range("B10").Interior.Color = 16777215
$ws.Range("B10").HorizontalAlignment = -4108

# Style E: A8/G8 reuse the plain/no-fill look already used elsewhere in the
# sheet (same as I2) -- copy that formatting across instead of inventing a
# new style.
$ws.Range("I2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("G8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Range("A8").Value = "Carrera 1 #67-21"
$ws.Range("G8").Value = 3108546364
